$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# The WeekApr8 sheet gets the schedule filled in (data entered from the
# project definition/blackboard) and becomes the active/selected tab.
# The WeekApr1 sheet loses the "selected" tab flag.
# ---------------------------------------------------------------------

$wsApr1 = $wb.Worksheets.Item("WeekApr1")
$wsApr8 = $wb.Worksheets.Item("WeekApr8")

# Fill in the WeekApr8 schedule data.
$wsApr8.Range("A3").Value = "Work on test plan"
$wsApr8.Range("B3").Value = "Andrew"
$wsApr8.Range("C3").Value = 2
$wsApr8.Range("D3").Value = "4/13/2018"

$wsApr8.Range("A4").Value = "Update buisness model to include expanations and test cases"
$wsApr8.Range("B4").Value = "Andrew"
$wsApr8.Range("C4").Value = 1
$wsApr8.Range("D4").Value = "4/13/2018"

$wsApr8.Range("A6").Value = "Create design alternatives"
$wsApr8.Range("B6").Value = "Kevin"
$wsApr8.Range("C6").Value = 1.15
$wsApr8.Range("D3").Copy()
$wsApr8.Range("D6").PasteSpecial(-4122)
$wsApr8.Range("D6").Value = "4/13/2018"

$wsApr8.Range("A7").Value = "Create design alternitive document"
$wsApr8.Range("B7").Value = "Kevin"
$wsApr8.Range("C7").Value = 2
$wsApr8.Range("D3").Copy()
$wsApr8.Range("D7").PasteSpecial(-4122)
$wsApr8.Range("D7").Value = "4/13/2018"

$wsApr8.Range("A9").Value = "Create PlayGame.java"
$wsApr8.Range("B9").Value = "Zach"
$wsApr8.Range("C9").Value = 0.5
$wsApr8.Range("D3").Copy()
$wsApr8.Range("D9").PasteSpecial(-4122)
$wsApr8.Range("D9").Value = "5/1/2018"

$wsApr8.Range("A10").Value = "Create Hand.java"
$wsApr8.Range("B10").Value = "Zach"
$wsApr8.Range("C10").Value = 0.5
$wsApr8.Range("D3").Copy()
$wsApr8.Range("D10").PasteSpecial(-4122)
$wsApr8.Range("D10").Value = "5/1/2018"

$wsApr8.Range("A11").Value = "Create BuildPhase.java"
$wsApr8.Range("B11").Value = "Zach"
$wsApr8.Range("C11").Value = 1
$wsApr8.Range("D3").Copy()
$wsApr8.Range("D11").PasteSpecial(-4122)
$wsApr8.Range("D11").Value = "5/1/2018"

$wsApr8.Range("A12").Value = "Create SpacePhase.java"
$wsApr8.Range("B12").Value = "Zach"
$wsApr8.Range("C12").Value = 1
$wsApr8.Range("D3").Copy()
$wsApr8.Range("D12").PasteSpecial(-4122)
$wsApr8.Range("D12").Value = "5/1/2018"

$wsApr8.Range("A13").Value = "Create Scorecard.java"
$wsApr8.Range("B13").Value = "Zach"
$wsApr8.Range("C13").Value = 0.5
$wsApr8.Range("D3").Copy()
$wsApr8.Range("D13").PasteSpecial(-4122)
$wsApr8.Range("D13").Value = "5/1/2018"

# These three "section header" rows only get a label in column A, with
# columns B:D left (un-styled/empty).
$wsApr8.Range("A8").Value = "Create source code based off current UML spec"
$wsApr8.Range("B8:D8").Clear()
$wsApr8.Range("A8").ClearFormats()

$wsApr8.Range("A5").Value = "Create a design alternitive analysis document"
$wsApr8.Range("B5:D5").Clear()
$wsApr8.Range("A5").ClearFormats()

$wsApr8.Range("A2").Value = "Create a test plan document"
$wsApr8.Range("B2:D2").Clear()
$wsApr8.Range("A2").ClearFormats()

# WeekApr8 becomes the active tab / selected sheet; WeekApr1 is no longer
# the selected tab (its own selection, E8, is left untouched).
$wsApr8.Activate()
$wsApr8.Range("A7").Select()
